$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A new batch of experiment task-order sheets was generated (new run
# timestamps), and the sheet tabs were re-created in a different order.
# Old tab order: GNG_TO, NB_TO, RS_TO, TOL_TO, vSAT_TO
# New tab order: vSAT_TO, RS_TO, GNG_TO, NB_TO, TOL_TO
# ---------------------------------------------------------------------------

# --- 1. Update the stim-file lists on each task-order sheet (by its current
#        name, before renaming) -------------------------------------------

$gng = $wb.Worksheets.Item("GNG_TO-16512555040398927")
$gng.Range("B2").Value = "go_stims-16515889312738192.csv"
$gng.Range("B3").Value = "GNG_stims-16515889312976053.csv"
$gng.Range("B4").Value = "go_stims-16515889312985995.csv"
$gng.Range("B5").Value = "GNG_stims-16515889313123982.csv"

$nb = $wb.Worksheets.Item("NB_TO-16512555072088923")
$nb.Range("B2").Value = "ZB-match_2-16515889321494348.csv"
$nb.Range("B3").Value = "ZB-match_7-16515889321109223.csv"
$nb.Range("B4").Value = "TB-16515889327464545.csv"
$nb.Range("B5").Value = "ZB-match_1-16515889313326857.csv"
$nb.Range("B6").Value = "TB-1651588932843642.csv"
$nb.Range("B7").Value = "OB-1651588932375385.csv"
$nb.Range("B8").Value = "TB-16515889331385667.csv"
$nb.Range("B9").Value = "OB-16515889326506908.csv"
$nb.Range("B10").Value = "OB-1651588932687131.csv"

# RS_TO (eyes open / eyes closed) keeps the same cell values - no content change.

$tol = $wb.Worksheets.Item("TOL_TO-16512555072738988")
$tol.Range("B2").Value = "MM_stims-16515889331941006.csv"
$tol.Range("B3").Value = "ZM_stims-16515889331692379.csv"
$tol.Range("B4").Value = "MM_stims-16515889332102146.csv"
$tol.Range("B5").Value = "ZM_stims-16515889331941006.csv"
$tol.Range("B6").Value = "MM_stims-16515889332260518.csv"
$tol.Range("B7").Value = "ZM_stims-16515889332112284.csv"

$vsat = $wb.Worksheets.Item("vSAT_TO-16512555073508916")
$vsat.Range("B2").Value = "vSAT_stims-16515889312488773.csv"
$vsat.Range("B3").Value = "vSAT_stims-16515889312334793.csv"
$vsat.Range("B4").Value = "SAT_stims-1651588931218104.csv"
$vsat.Range("B5").Value = "SAT_stims-16515889312043374.csv"

# --- 2. Rename every sheet to its freshly generated tab name ---------------

$wb.Worksheets.Item("GNG_TO-16512555040398927").Name = "GNG_TO-16515889313135014"
$wb.Worksheets.Item("NB_TO-16512555072088923").Name = "NB_TO-16515889331621997"
$wb.Worksheets.Item("RS_TO-16512555072148976").Name = "RS_TO-16515889312700062"
$wb.Worksheets.Item("TOL_TO-16512555072738988").Name = "TOL_TO-16515889332270248"
$wb.Worksheets.Item("vSAT_TO-16512555073508916").Name = "vSAT_TO-16515889312651088"

# --- 3. Re-order the tabs: vSAT_TO, RS_TO, GNG_TO, NB_TO, TOL_TO -----------
# Re-fetch each sheet right before moving it so the reference is never stale.

$wb.Worksheets.Item("vSAT_TO-16515889312651088").Move($wb.Worksheets.Item("GNG_TO-16515889313135014"))
$wb.Worksheets.Item("RS_TO-16515889312700062").Move($wb.Worksheets.Item("GNG_TO-16515889313135014"))
